# Applies the commit's edits: splits several runs so that the words Word's
# spell-checker flags ("subthema", "Mujagic", "Macroeconomische", "graaiflatie",
# "krimpflatie", "notules", "jammergenoeg", "foodwatch") sit in their own run,
# bracketed by <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>
# markers - exactly like Word does automatically once it has proofed the text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate the paragraph that contains $anchorText, replace its content
# (everything except the paragraph mark, so pPr / paraId / rsids survive)
# with the supplied OOXML run fragment.
# ---------------------------------------------------------------------------
function Replace-ParagraphContent {
    param(
        [string]$anchorText,
        [string]$innerXml
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor not found: $anchorText"
    }

    $rng.Expand(4) | Out-Null   # wdParagraph
    $t = $rng.Text
    $lastChar = $t.Substring($t.Length - 1)
    if ($lastChar -eq [char]13) {
        $rng.End = $rng.End - 1
    }
    $rng.Delete()

    $xmlHead = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
    $xmlTail = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $xml = $xmlHead + $innerXml + $xmlTail
    $rng.InsertXML($xml)
}

function MakeRun {
    param(
        [string]$rpr,
        [string]$text,
        [bool]$preserve
    )
    $sp = ""
    if ($preserve) { $sp = ' xml:space="preserve"' }
    $result = '<w:r><w:rPr>' + $rpr + '</w:rPr><w:t' + $sp + '>' + $text + '</w:t></w:r>'
    return $result
}

$rprBold  = '<w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:b/><w:bCs/><w:lang w:val="nl-NL"/>'
$rprPlain = '<w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:lang w:val="nl-NL"/>'
$spellStart = '<w:proofErr w:type="spellStart"/>'
$spellEnd   = '<w:proofErr w:type="spellEnd"/>'

# ---------------------------------------------------------------------------
# 1) "Alinea 2: (subthema) (ECB)" -> split the " (subthema)" run
# ---------------------------------------------------------------------------
$p1a = MakeRun $rprBold 'Alinea 2:' $false
$p1b = MakeRun $rprBold ' (' $true
$p1c = MakeRun $rprBold 'subthema' $false
$p1d = MakeRun $rprBold ')' $false
$p1e = MakeRun $rprBold ' (ECB)' $true
$inner1 = $p1a + $p1b + $spellStart + $p1c + $spellEnd + $p1d + $p1e
Replace-ParagraphContent 'Alinea 2: (subthema) (ECB)' $inner1

# ---------------------------------------------------------------------------
# 2) "Dat wil ik bespreken met econoom Edin Mujagic wie voor verschillende
#    instanties het Macroeconomische nieuws brengt."
# ---------------------------------------------------------------------------
$p2a = MakeRun $rprPlain 'Dat wil ik bespreken met econoom Edin ' $true
$p2b = MakeRun $rprPlain 'Mujagic' $false
$p2c = MakeRun $rprPlain ' wie voor verschillende instanties het ' $true
$p2d = MakeRun $rprPlain 'Macroeconomische' $false
$p2e = MakeRun $rprPlain ' nieuws brengt.' $true
$inner2 = $p2a + $spellStart + $p2b + $spellEnd + $p2c + $spellStart + $p2d + $spellEnd + $p2e
Replace-ParagraphContent 'Dat wil ik bespreken met econoom Edin Mujagic' $inner2

# ---------------------------------------------------------------------------
# 3) "Alinea 3: (subthema)" -> split the " (subthema)" run
# ---------------------------------------------------------------------------
$p3a = MakeRun $rprBold 'Alinea 3:' $false
$p3b = MakeRun $rprBold ' (' $true
$p3c = MakeRun $rprBold 'subthema' $false
$p3d = MakeRun $rprBold ')' $false
$inner3 = $p3a + $p3b + $spellStart + $p3c + $spellEnd + $p3d
Replace-ParagraphContent 'Alinea 3: (subthema)' $inner3

# ---------------------------------------------------------------------------
# 4) "Resultaat van beleid ECB en krimp- & graaiflatie, ... 'graaiflatie' en
#    'krimpflatie'. Hoe kunnen we dit herkennen ..."
# ---------------------------------------------------------------------------
$lsq = [char]0x2018
$rsq = [char]0x2019
$p4a = MakeRun $rprPlain 'Resultaat van beleid ECB en krimp- &amp; ' $true
$p4b = MakeRun $rprPlain 'graaiflatie' $false
$p4cText = ', deze punten hebben geleid tot de huidige situatie voor het afgelopen jaar. De prijzen van consumentengoederen daalt maar niet en er treden nu ook nieuwe economische fenomenen op zoals ' + $lsq
$p4c = MakeRun $rprPlain $p4cText $false
$p4d = MakeRun $rprPlain 'graaiflatie' $false
$p4eText = $rsq + ' en ' + $lsq
$p4e = MakeRun $rprPlain $p4eText $false
$p4f = MakeRun $rprPlain 'krimpflatie' $false
$p4gText = $rsq + '. Hoe kunnen we dit herkennen en hoe gaat dit in de toekomst eruit zien?'
$p4g = MakeRun $rprPlain $p4gText $false
$inner4 = $p4a + $spellStart + $p4b + $spellEnd + $p4c + $spellStart + $p4d + $spellEnd + $p4e + $spellStart + $p4f + $spellEnd + $p4g
Replace-ParagraphContent 'Resultaat van beleid ECB' $inner4

# ---------------------------------------------------------------------------
# 5) "Na contact met de ECB word ik verwezen naar de meest recente uitspraken
#    en notules van hun president Christine Lagarde. ..."
# ---------------------------------------------------------------------------
$p5a = MakeRun $rprPlain 'Na contact met de ECB word ik verwezen naar de meest recente uitspraken en ' $true
$p5b = MakeRun $rprPlain 'notules' $false
$p5c = MakeRun $rprPlain ' van hun president Christine Lagarde. Mijn vraag of de ECB nog steeds voet bij stuk houd word indirect beantwoord met het volgende: hun huidige doelstelling is bekend gemaakt tijdens hun persconferentie van 26 November 2023, waarin besproken word in 2025 weer een rentestand te hebben van 2%. ' $true
$inner5 = $p5a + $spellStart + $p5b + $spellEnd + $p5c
Replace-ParagraphContent 'Na contact met de ECB' $inner5

# ---------------------------------------------------------------------------
# 6) "Supermarktketens weigeren jammergenoeg mijn vragen over verandering qua
#    prijs in de supermarkt, maar in reactie naar een [voedsel autoriteit]
#    leggen ze de schuld ..." - only the first run is split, the highlighted
#    run and the trailing run are untouched.
# ---------------------------------------------------------------------------
$p6a = MakeRun $rprPlain 'Supermarktketens weigeren ' $true
$p6b = MakeRun $rprPlain 'jammergenoeg' $false
$p6c = MakeRun $rprPlain ' mijn vragen over verandering qua prijs in de supermarkt, maar in reactie naar een ' $true
$rprHighlight = '<w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/><w:highlight w:val="yellow"/><w:lang w:val="nl-NL"/>'
$p6d = MakeRun $rprHighlight 'voedsel autoriteit' $false
$p6e = MakeRun $rprPlain ' leggen ze de schuld van deze duurdere producten bij de fabrikant. Ook word er verwezen naar het CBL de branchevereniging voor levensmiddelen.' $true
$inner6 = $p6a + $spellStart + $p6b + $spellEnd + $p6c + $p6d + $p6e
Replace-ParagraphContent 'Supermarktketens weigeren jammergenoeg' $inner6

# ---------------------------------------------------------------------------
# 7) "Initiatieven zoals foodwatch bestuderen deze producten en tonen zowel
#    prijsstijgingen als inhoudsafnames"
# ---------------------------------------------------------------------------
$p7a = MakeRun $rprPlain 'Initiatieven zoals ' $true
$p7b = MakeRun $rprPlain 'foodwatch' $false
$p7c = MakeRun $rprPlain ' bestuderen deze producten en tonen zowel prijsstijgingen als inhoudsafnames' $true
$inner7 = $p7a + $spellStart + $p7b + $spellEnd + $p7c
Replace-ParagraphContent 'Initiatieven zoals foodwatch' $inner7

Write-Host "Done applying edits"
